$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like literal string into a cell without Excel
# auto-converting it to a date serial number. We build the text via a
# formula (so it is never typed as a raw literal), then copy/paste-special
# the computed value back over itself (xlPasteValues = -4163). This keeps
# the cell as a plain shared-string text cell with the default style,
# matching how the source data file stored these dates.
function Set-TextDate($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 53
Set-TextDate $ws.Cells.Item(53, 1) "2018.08.30"
$ws.Cells.Item(53, 2).Value = "15:39:37"
$ws.Cells.Item(53, 3).Value = "RS"
$ws.Cells.Item(53, 4).Value = 10
$ws.Cells.Item(53, 5).Value = 250
$ws.Cells.Item(53, 6).Value = "N/A"
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 681
$ws.Cells.Item(53, 9).Value = 0.0316
$ws.Cells.Item(53, 10).Value = 4
$ws.Cells.Item(53, 11).Value = 35.30099402320572
$ws.Cells.Item(53, 12).Value = 0.01523220108113661

# Row 54
Set-TextDate $ws.Cells.Item(54, 1) "2018.08.30"
$ws.Cells.Item(54, 2).Value = "15:39:38"
$ws.Cells.Item(54, 3).Value = "RS"
$ws.Cells.Item(54, 4).Value = 10
$ws.Cells.Item(54, 5).Value = 250
$ws.Cells.Item(54, 6).Value = 0.1
$ws.Cells.Item(54, 7).Value = 0.1
$ws.Cells.Item(54, 8).Value = 254
$ws.Cells.Item(54, 9).Value = 0.0216
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = "N/A"
$ws.Cells.Item(54, 12).Value = "N/A"

# Row 55
Set-TextDate $ws.Cells.Item(55, 1) "2018.08.30"
$ws.Cells.Item(55, 2).Value = "15:39:40"
$ws.Cells.Item(55, 3).Value = "RS"
$ws.Cells.Item(55, 4).Value = 10
$ws.Cells.Item(55, 5).Value = 250
$ws.Cells.Item(55, 6).Value = 0.1
$ws.Cells.Item(55, 7).Value = 0.25
$ws.Cells.Item(55, 8).Value = 343
$ws.Cells.Item(55, 9).Value = 0.0236
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = "N/A"
$ws.Cells.Item(55, 12).Value = "N/A"

# Row 56
Set-TextDate $ws.Cells.Item(56, 1) "2018.08.30"
$ws.Cells.Item(56, 2).Value = "15:39:42"
$ws.Cells.Item(56, 3).Value = "RS"
$ws.Cells.Item(56, 4).Value = 10
$ws.Cells.Item(56, 5).Value = 250
$ws.Cells.Item(56, 6).Value = 0.1
$ws.Cells.Item(56, 7).Value = 0.75
$ws.Cells.Item(56, 8).Value = 725
$ws.Cells.Item(56, 9).Value = 0.0286
$ws.Cells.Item(56, 10).Value = 2
$ws.Cells.Item(56, 11).Value = 49.54912798466926
$ws.Cells.Item(56, 12).Value = 0.01021924049293189

# Row 57
Set-TextDate $ws.Cells.Item(57, 1) "2018.08.30"
$ws.Cells.Item(57, 2).Value = "15:41:38"
$ws.Cells.Item(57, 3).Value = "RS"
$ws.Cells.Item(57, 4).Value = 10
$ws.Cells.Item(57, 5).Value = 250
$ws.Cells.Item(57, 6).Value = "N/A"
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 681
$ws.Cells.Item(57, 9).Value = 0.0266
$ws.Cells.Item(57, 10).Value = 4
$ws.Cells.Item(57, 11).Value = 35.30099402320572
$ws.Cells.Item(57, 12).Value = 0.01523220108113661

# Row 58
Set-TextDate $ws.Cells.Item(58, 1) "2018.08.30"
$ws.Cells.Item(58, 2).Value = "15:43:46"
$ws.Cells.Item(58, 3).Value = "RS"
$ws.Cells.Item(58, 4).Value = 10
$ws.Cells.Item(58, 5).Value = 250
$ws.Cells.Item(58, 6).Value = "N/A"
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 681
$ws.Cells.Item(58, 9).Value = 0.0276
$ws.Cells.Item(58, 10).Value = 4
$ws.Cells.Item(58, 11).Value = 35.30099402320572
$ws.Cells.Item(58, 12).Value = 0.01523220108113661

# Row 59
Set-TextDate $ws.Cells.Item(59, 1) "2018.08.30"
$ws.Cells.Item(59, 2).Value = "16:17:22"
$ws.Cells.Item(59, 3).Value = "RS"
$ws.Cells.Item(59, 4).Value = 10
$ws.Cells.Item(59, 5).Value = 250
$ws.Cells.Item(59, 6).Value = "N/A"
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 681
$ws.Cells.Item(59, 9).Value = 0.0276
$ws.Cells.Item(59, 10).Value = 4
$ws.Cells.Item(59, 11).Value = 35.30099402320572
$ws.Cells.Item(59, 12).Value = 0.01523220108113661

# Row 60
Set-TextDate $ws.Cells.Item(60, 1) "2018.08.30"
$ws.Cells.Item(60, 2).Value = "16:17:24"
$ws.Cells.Item(60, 3).Value = "RS"
$ws.Cells.Item(60, 4).Value = -2
$ws.Cells.Item(60, 5).Value = 250
$ws.Cells.Item(60, 6).Value = "N/A"
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 45
$ws.Cells.Item(60, 9).Value = 0.0171
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = "N/A"
$ws.Cells.Item(60, 12).Value = "N/A"

# Row 61
Set-TextDate $ws.Cells.Item(61, 1) "2018.08.30"
$ws.Cells.Item(61, 2).Value = "16:17:25"
$ws.Cells.Item(61, 3).Value = "RS"
$ws.Cells.Item(61, 4).Value = 10
$ws.Cells.Item(61, 5).Value = 250
$ws.Cells.Item(61, 6).Value = "N/A"
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 681
$ws.Cells.Item(61, 9).Value = 0.0236
$ws.Cells.Item(61, 10).Value = 4
$ws.Cells.Item(61, 11).Value = 35.30099402320572
$ws.Cells.Item(61, 12).Value = 0.01523220108113661

# Row 62
Set-TextDate $ws.Cells.Item(62, 1) "2018.08.30"
$ws.Cells.Item(62, 2).Value = "16:17:26"
$ws.Cells.Item(62, 3).Value = "RS"
$ws.Cells.Item(62, 4).Value = 10
$ws.Cells.Item(62, 5).Value = 250
$ws.Cells.Item(62, 6).Value = 0.1
$ws.Cells.Item(62, 7).Value = 0.5
$ws.Cells.Item(62, 8).Value = 544
$ws.Cells.Item(62, 9).Value = 0.0256
$ws.Cells.Item(62, 10).Value = 1
$ws.Cells.Item(62, 11).Value = 112.622406108482
$ws.Cells.Item(62, 12).Value = "N/A"

# Row 63
Set-TextDate $ws.Cells.Item(63, 1) "2018.08.30"
$ws.Cells.Item(63, 2).Value = "16:20:10"
$ws.Cells.Item(63, 3).Value = "RS"
$ws.Cells.Item(63, 4).Value = 10
$ws.Cells.Item(63, 5).Value = 250
$ws.Cells.Item(63, 6).Value = "N/A"
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 681
$ws.Cells.Item(63, 9).Value = 0.0261
$ws.Cells.Item(63, 10).Value = 4
$ws.Cells.Item(63, 11).Value = 35.30099402320572
$ws.Cells.Item(63, 12).Value = 0.01523220108113661

# Row 64
Set-TextDate $ws.Cells.Item(64, 1) "2018.08.30"
$ws.Cells.Item(64, 2).Value = "17:00:29"
$ws.Cells.Item(64, 3).Value = "RS"
$ws.Cells.Item(64, 4).Value = 10
$ws.Cells.Item(64, 5).Value = 250
$ws.Cells.Item(64, 6).Value = "N/A"
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 681
$ws.Cells.Item(64, 9).Value = 0.0248
$ws.Cells.Item(64, 10).Value = 4
$ws.Cells.Item(64, 11).Value = 35.30099402320572
$ws.Cells.Item(64, 12).Value = 0.01523220108113661

# Row 65
Set-TextDate $ws.Cells.Item(65, 1) "2018.08.30"
$ws.Cells.Item(65, 2).Value = "17:03:39"
$ws.Cells.Item(65, 3).Value = "RS"
$ws.Cells.Item(65, 4).Value = 10
$ws.Cells.Item(65, 5).Value = 250
$ws.Cells.Item(65, 6).Value = "N/A"
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 681
$ws.Cells.Item(65, 9).Value = 0.0246
$ws.Cells.Item(65, 10).Value = 4
$ws.Cells.Item(65, 11).Value = 35.30099402320572
$ws.Cells.Item(65, 12).Value = 0.01523220108113661

# Row 66
Set-TextDate $ws.Cells.Item(66, 1) "2018.08.30"
$ws.Cells.Item(66, 2).Value = "17:03:42"
$ws.Cells.Item(66, 3).Value = "RS"
$ws.Cells.Item(66, 4).Value = 10
$ws.Cells.Item(66, 5).Value = 250
$ws.Cells.Item(66, 6).Value = 0.1
$ws.Cells.Item(66, 7).Value = 0.97
$ws.Cells.Item(66, 8).Value = 793
$ws.Cells.Item(66, 9).Value = 0.0221
$ws.Cells.Item(66, 10).Value = 3
$ws.Cells.Item(66, 11).Value = 36.523440888569
$ws.Cells.Item(66, 12).Value = 0.01617273454703009

# Row 67
Set-TextDate $ws.Cells.Item(67, 1) "2018.08.30"
$ws.Cells.Item(67, 2).Value = "17:03:46"
$ws.Cells.Item(67, 3).Value = "RS"
$ws.Cells.Item(67, 4).Value = 10
$ws.Cells.Item(67, 5).Value = 250
$ws.Cells.Item(67, 6).Value = 0.1
$ws.Cells.Item(67, 7).Value = 0.93
$ws.Cells.Item(67, 8).Value = 814
$ws.Cells.Item(67, 9).Value = 0.0271
$ws.Cells.Item(67, 10).Value = 3
$ws.Cells.Item(67, 11).Value = 38.21998138918514
$ws.Cells.Item(67, 12).Value = 0.01462097067859595

